# Update "想去人数" (column F) counts on the 展览 (Expo), 演出 (Show) and
# 全部类型 (All Types) sheets to reflect the refreshed gh-pages data pull.
# (本地生活 sheet is unchanged in this refresh.)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1296
$wsExpo.Range("F6").Value = 384
$wsExpo.Range("F7").Value = 1211
$wsExpo.Range("F8").Value = 454
$wsExpo.Range("F9").Value = 7409
$wsExpo.Range("F11").Value = 97
$wsExpo.Range("F12").Value = 2061
$wsExpo.Range("F13").Value = 8064
$wsExpo.Range("F16").Value = 5542
$wsExpo.Range("F17").Value = 57
$wsExpo.Range("F18").Value = 2469
$wsExpo.Range("F19").Value = 1048
$wsExpo.Range("F20").Value = 4567
$wsExpo.Range("F22").Value = 389
$wsExpo.Range("F24").Value = 15
$wsExpo.Range("F25").Value = 411
$wsExpo.Range("F26").Value = 655
$wsExpo.Range("F28").Value = 2502
$wsExpo.Range("F30").Value = 279
$wsExpo.Range("F31").Value = 93
$wsExpo.Range("F32").Value = 185
$wsExpo.Range("F33").Value = 608
$wsExpo.Range("F36").Value = 1555
$wsExpo.Range("F38").Value = 11
$wsExpo.Range("F39").Value = 2445
$wsExpo.Range("F40").Value = 2234
$wsExpo.Range("F43").Value = 17

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 23

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 384
$wsAll.Range("F7").Value = 1211
$wsAll.Range("F8").Value = 454
$wsAll.Range("F9").Value = 7409
$wsAll.Range("F11").Value = 97
$wsAll.Range("F12").Value = 2061
$wsAll.Range("F13").Value = 8064
$wsAll.Range("F16").Value = 5542
$wsAll.Range("F17").Value = 57
$wsAll.Range("F18").Value = 2469
$wsAll.Range("F19").Value = 1048
$wsAll.Range("F20").Value = 4567
$wsAll.Range("F22").Value = 389
$wsAll.Range("F25").Value = 15
$wsAll.Range("F27").Value = 411
$wsAll.Range("F28").Value = 655
$wsAll.Range("F30").Value = 2502
$wsAll.Range("F32").Value = 279
$wsAll.Range("F33").Value = 93
$wsAll.Range("F34").Value = 185
$wsAll.Range("F35").Value = 23
$wsAll.Range("F36").Value = 608
$wsAll.Range("F40").Value = 1555
$wsAll.Range("F42").Value = 11
$wsAll.Range("F43").Value = 2445
$wsAll.Range("F45").Value = 2234
$wsAll.Range("F48").Value = 17
